$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 to hold the Portuguese ("por") language record
$ws.Range("A3").Value = "por"
$ws.Range("B3").Value = "Português"
$ws.Range("C3").Value = "Português"
$ws.Range("D3").Value = "Portuguese"

# Remove the remaining language rows (old rows 4-9: ara, kan, hin, tam, por, asdf)
$ws.Range("A4:E9").EntireRow.Delete() | Out-Null

# Match the saved selection state from the authored workbook
$ws.Range("C3").Select() | Out-Null
